# Apply the "break out stock.yaml completed" edit to the "10per change" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Fix E53 and E54: they were stored as inline strings ("543287", "590024")
# but should be numeric values.
$ws.Range("E53").Value = 543287
$ws.Range("E54").Value = 590024

# Append new row 55 with the latest FACT screener entry.
$ws.Range("A55").Value = "25/06/2024 09:44:40"
$ws.Range("B55").Value = 1
$ws.Range("C55").Value = "FACT"
$ws.Range("D55").Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Range("E55").NumberFormat = "@"
$ws.Range("E55").Value = "590024"
$ws.Range("E55").Style = "Normal"
$ws.Range("F55").Value = -1.97
$ws.Range("G55").Value = 1000.85
$ws.Range("H55").Value = 1906971
